$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$caseName = $ws.Range("B6").Value2
$inputFile = $ws.Range("C6").Value2
$outputFile = $ws.Range("D6").Value2
$initPF = $ws.Range("F6").Value2

$ws.Range("A7").Value = 33
$ws.Range("B7").Value = $caseName
$ws.Range("C7").Value = $inputFile
$ws.Range("D7").Value = $outputFile
$ws.Range("E7").Value = "2024-02-16T04:07:11.353"
$ws.Range("F7").Value = $initPF
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = $false
$ws.Range("I7").Value = $false

$ws.Range("A8").Value = 33
$ws.Range("B8").Value = $caseName
$ws.Range("C8").Value = $inputFile
$ws.Range("D8").Value = $outputFile
$ws.Range("E8").Value = "2024-02-16T04:08:15.245"
$ws.Range("F8").Value = $initPF
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false
$ws.Range("I8").Value = $false
